$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Frage 1")

# --- New column F values: "Antragsanzahl aus Parlamentsbeantwortung 18726/AB" ---
# Header for column F (row 4)
$ws.Range("F4").Value = "Antragsanzahl aus Parlamentsbeantwortung 18726/AB (nicht für alle Fachgebiete verfügbar)"

# Data values for column F (only some rows have data available)
$ws.Range("F5").Value = 131668
$ws.Range("F7").Value = 80808
$ws.Range("F8").Value = 20518
$ws.Range("F9").Value = 80035
$ws.Range("F10").Value = 87008
$ws.Range("F12").Value = 51537
$ws.Range("F16").Value = 77611
$ws.Range("F20").Value = 32213
$ws.Range("F24").Value = 23418
$ws.Range("F34").Value = 1676

# Column G: ratio of F (Parlament answer) to C (SVS count) for every data row 5..34
$ws.Range("G5").Formula = "=F5/C5"
$ws.Range("G6:G34").Formula = "=F6/C6"

# Row 43: headers for the two summary columns
$ws.Range("G43").Value = "Summe für Fachgebiete wo Daten verfügbar"
$ws.Range("H43").Value = "Summe für Fachgebiete wo keine Daten verfügbar"
$ws.Range("G43:H43").Font.Bold = $true

# Row 44: totals
$ws.Range("F44").Formula = "=SUM(F5:F41)"
$ws.Range("G44").Formula = '=SUMIF(F5:F41,"<>",C5:C41)'
$ws.Range("H44").Formula = '=SUMIF(F5:F41,"",C5:C41)'

# Row 46/47: Zahnarztanträge = Zahn-,Mund- und Kieferheilkunde (C21) + Zahnarzt (C30)
$ws.Range("G46").Value = "Zahnarztanträge"
$ws.Range("G47").Formula = "=C30+C21"

# Selection moved to B34 on the active sheet
$ws.Activate()
$ws.Range("B34").Select() | Out-Null
